$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New strikeout (K) values regenerated from source data, replacing the
# previous "Strike#" derived values in column G for rows 2-36.
$kValues = @{
    2  = 4
    3  = 2
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 2
    9  = 5
    10 = 4
    11 = 5
    12 = 9
    13 = 4
    14 = 8
    15 = 7
    16 = 4
    17 = 5
    18 = 4
    19 = 10
    20 = 6
    21 = 8
    22 = 7
    23 = 5
    24 = 1
    25 = 5
    26 = 10
    27 = 5
    28 = 10
    29 = 6
    30 = 11
    31 = 4
    32 = 5
    33 = 6
    34 = 8
    35 = 1
    36 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
